$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, pushing the existing rows 74-193 down to 75-194.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with its data (same shape as the other
# "Rabanito" rows in the block, with its own date / price figures).
$ws.Cells.Item(74, 1).Value = 9
$ws.Cells.Item(74, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74, 3).Value = "Metropolitana"
$ws.Cells.Item(74, 4).Value = [DateTime]"2021-12-03"
$ws.Cells.Item(74, 5).Value = 13
$ws.Cells.Item(74, 6).Value = 300000001
$ws.Cells.Item(74, 7).Value = "Rabanito"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 7900
$ws.Cells.Item(74, 11).Value = 2500
$ws.Cells.Item(74, 12).Value = 3000
$ws.Cells.Item(74, 13).Value = 2747
$ws.Cells.Item(74, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(74, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(74, 16).Value = 27
$ws.Cells.Item(74, 17).Value = 100
$ws.Cells.Item(74, 18).Value = "Hortaliza"
